$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (239-244): date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44313, 0, 7, 152.6717557251908),
    @(44314, 0, 7, 152.6717557251908),
    @(44315, 0, 4, 87.24100327153762),
    @(44316, 1, 1, 21.81025081788441),
    @(44317, 0, 1, 21.81025081788441),
    @(44318, 0, 1, 21.81025081788441)
)

$startRow = 239
$endRow = 244

# Copy formatting from the last existing data row (238) down to the new rows
$ws.Range("A238:D238").Copy() | Out-Null
$ws.Range("A239:D244").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
